$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.727.34"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.333.16"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Formula = "=""578.45"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Formula = "=""174.94"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.329.70"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Formula = "=""0.581"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Formula = "=""46.26"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Formula = "=""0.0000272"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Formula = "=""706.03"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "3.871.21"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "67.756.86"
$ws.Range("D19").Value = "3.339.75"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Formula = "=""17.38"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Formula = "=""11.00"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").Formula = "=""0.894"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Formula = "=""5.39"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").Formula = "=""16.94"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Formula = "=""98.56"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Formula = "=""2.69"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").Formula = "=""9.42"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Formula = "=""33.23"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("E31").Value = "  +5.38%  "
$ws.Range("D32").Formula = "=""569.65"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").Formula = "=""10.98"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "3.704.87"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("D37").Formula = "=""56.97"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +2.79%  "
$ws.Range("D38").Formula = "=""3.32"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Formula = "=""34.18"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +6.40%  "
$ws.Range("D40").Formula = "=""0.131"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Formula = "=""2.65"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Formula = "=""3.17"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "0.0₃0675"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Formula = "=""0.337"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("E47").Value = "  +6.39%  "
$ws.Range("D48").Formula = "=""0.128"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  -5.04%  "
$ws.Range("D51").Formula = "=""129.06"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.13%  "
$excel.CutCopyMode = 0
